$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G2: total hits for "digital tvilling" updated from 22 to 26
$ws.Range("G2").Value = 26

# Row 12: " bim " search term updated to include hyphen variant " bim[- ]"
$ws.Range("A12").Value = " bim[- ]"
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = 4
$ws.Range("E12").Value = 6
